$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '64.287.33'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').Value = '3.507.94'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue 'D5' '584.69'
$ws.Range('E5').Value = '  -0.40%  '
Set-TextValue 'D6' '134.86'
$ws.Range('E6').Value = '  +1.14%  '
$ws.Range('D7').Value = '3.508.25'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -0.14%  '
Set-TextValue 'D10' '0.125'
$ws.Range('E10').Value = '  +0.24%  '
Set-TextValue 'D11' '7.11'
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('E12').Value = '  -2.17%  '
$ws.Range('D13').Value = '4.104.96'
$ws.Range('E13').Value = '  -0.06%  '
Set-TextValue 'D14' '27.41'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('E16').Value = '  +1.16%  '
$ws.Range('D17').Value = '3.506.36'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '64.291.19'
$ws.Range('E18').Value = '  -0.78%  '
Set-TextValue 'D19' '9.82'
$ws.Range('E19').Value = '  -1.75%  '
Set-TextValue 'D20' '13.90'
$ws.Range('E20').Value = '  -2.49%  '
Set-TextValue 'D21' '5.60'
$ws.Range('E21').Value = '  -1.07%  '
Set-TextValue 'D22' '384.37'
$ws.Range('E22').Value = '  -1.63%  '
Set-TextValue 'D23' '0.569'
$ws.Range('E23').Value = '  -1.25%  '
$ws.Range('D24').Value = '3.648.82'
$ws.Range('E24').Value = '  -0.14%  '
Set-TextValue 'D25' '73.93'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('E26').Value = '  +0.00%  '
Set-TextValue 'D27' '5.72'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  +5.47%  '
$ws.Range('E29').Value = '  -0.36%  '
Set-TextValue 'D30' '7.62'
$ws.Range('E30').Value = '  +2.12%  '
$ws.Range('E31').Value = '  -0.01%  '
Set-TextValue 'D32' '8.36'
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('E33').Value = '  -2.00%  '
$ws.Range('D34').Value = '3.520.45'
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('E36').Value = '  +0.34%  '
Set-TextValue 'D37' '23.62'
$ws.Range('E37').Value = '  -1.45%  '
Set-TextValue 'D38' '5.34'
$ws.Range('E38').Value = '  +3.03%  '
Set-TextValue 'D39' '1.56'
$ws.Range('E39').Value = '  -2.30%  '
Set-TextValue 'D40' '6.90'
$ws.Range('E40').Value = '  -0.82%  '
Set-TextValue 'D41' '163.99'
$ws.Range('E41').Value = '  -4.56%  '
Set-TextValue 'D42' '0.0786'
$ws.Range('E42').Value = '  -2.75%  '
$ws.Range('E43').Value = '  -0.73%  '
Set-TextValue 'D44' '26.20'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E46').Value = '  -0.49%  '
Set-TextValue 'D47' '41.81'
$ws.Range('E47').Value = '  -1.06%  '
Set-TextValue 'D48' '4.40'
$ws.Range('E48').Value = '  -0.05%  '
Set-TextValue 'D49' '1.64'
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('D50').Value = '2.483.22'
$ws.Range('E50').Value = '  -0.14%  '
Set-TextValue 'D51' '6.77'
$ws.Range('E51').Value = '  -1.38%  '
